# Add files via upload
# Update the Approved/Rejected decision (col I) and add a note (col J)
# for the two rows that were reviewed ("New Case" scenario rows 12 & 20
# of the Test-Cases table), then leave the selection where the user left
# it after making the edits (cell J10, scrolled so column H is leftmost).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Activate()

$ws.Range("I12").Value = "Rejected"
$ws.Range("J12").Value = "test"

$ws.Range("I20").Value = "Rejected"
$ws.Range("J20").Value = "test"

# Move the view/selection to match where the author ended up.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
[void]$ws.Range("J10").Select()
